$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 24
$ws.Range("E2").Value = 31
$ws.Range("F2").Value = 77.41935483870968
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 22.58064516129032
$ws.Range("E4").Value = 31
$ws.Range("E5").Value = 31
$ws.Range("E6").Value = 31
$ws.Range("E7").Value = 31
$ws.Range("E8").Value = 31
$ws.Range("E9").Value = 31
$ws.Range("E10").Value = 31
$ws.Range("E11").Value = 31
$ws.Range("E12").Value = 31
$ws.Range("E13").Value = 31
$ws.Range("D14").Value = 32
$ws.Range("E14").Value = 34
$ws.Range("F14").Value = 94.11764705882352
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 34
$ws.Range("F15").Value = 5.88235294117647
$ws.Range("E16").Value = 34
$ws.Range("E17").Value = 34
$ws.Range("E18").Value = 34
$ws.Range("E19").Value = 34
$ws.Range("E20").Value = 34
$ws.Range("E21").Value = 34
$ws.Range("E22").Value = 34
$ws.Range("E23").Value = 34
$ws.Range("E24").Value = 34
$ws.Range("E25").Value = 34
$ws.Range("D26").Value = 28
$ws.Range("E26").Value = 30
$ws.Range("F26").Value = 93.33333333333333
$ws.Range("E27").Value = 30
$ws.Range("F27").Value = 6.666666666666667
$ws.Range("E28").Value = 30
$ws.Range("E29").Value = 30
$ws.Range("E30").Value = 30
$ws.Range("E31").Value = 30
$ws.Range("E32").Value = 30
$ws.Range("E33").Value = 30
$ws.Range("E34").Value = 30
$ws.Range("E35").Value = 30
$ws.Range("E36").Value = 30
$ws.Range("E37").Value = 30
$ws.Range("D38").Value = 13
$ws.Range("E38").Value = 18
$ws.Range("F38").Value = 72.22222222222221
$ws.Range("D39").Value = 3
$ws.Range("E39").Value = 18
$ws.Range("F39").Value = 16.66666666666666
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 18
$ws.Range("F40").Value = 5.555555555555555
$ws.Range("A41").Value = 39
$ws.Range("C41").Value = "Ocean Alkalinisation"
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 18
$ws.Range("F41").Value = 5.555555555555555
$ws.Range("A42").Value = 38
$ws.Range("C42").Value = "Ocean fertilisation"
$ws.Range("E42").Value = 18
$ws.Range("E43").Value = 18
$ws.Range("E44").Value = 18
$ws.Range("E45").Value = 18
$ws.Range("E46").Value = 18
$ws.Range("E47").Value = 18
$ws.Range("E48").Value = 18
$ws.Range("E49").Value = 18
$ws.Range("D50").Value = 7
$ws.Range("E50").Value = 20
$ws.Range("F50").Value = 35
$ws.Range("E51").Value = 20
$ws.Range("F51").Value = 25
$ws.Range("D52").Value = 5
$ws.Range("E52").Value = 20
$ws.Range("F52").Value = 25
$ws.Range("D53").Value = 3
$ws.Range("E53").Value = 20
$ws.Range("F53").Value = 15
$ws.Range("E54").Value = 20
$ws.Range("E55").Value = 20
$ws.Range("E56").Value = 20
$ws.Range("E57").Value = 20
$ws.Range("E58").Value = 20
$ws.Range("E59").Value = 20
$ws.Range("E60").Value = 20
$ws.Range("E61").Value = 20
$ws.Range("D62").Value = 17
$ws.Range("E62").Value = 17
$ws.Range("E63").Value = 17
$ws.Range("E64").Value = 17
$ws.Range("E65").Value = 17
$ws.Range("E66").Value = 17
$ws.Range("E67").Value = 17
$ws.Range("E68").Value = 17
$ws.Range("E69").Value = 17
$ws.Range("E70").Value = 17
$ws.Range("E71").Value = 17
$ws.Range("E72").Value = 17
$ws.Range("E73").Value = 17
